$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume refresh (GitHub Actions scheduled update).
# For each changed row, set the Price (D) and Volume(1h) (E) columns.
# The source cells are plain text (inlineStr) containing numeric-looking
# strings (e.g. "328.29") and percentages (e.g. "-0.18%"), so we must
# force Excel to keep them as text (quote-prefix) rather than silently
# re-typing them as numbers, then reset the style so no extra cell
# formatting is introduced.

$ws.Range("D2").Value = "'328.29"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.18%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'44.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-0.02%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.070"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-7.63%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08392"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'4.12%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.947"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-5.14%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.9755"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.12%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'2.501"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-4.58%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1143"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.08%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1900"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.47%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09655"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-2.77%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04636"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.64%"
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "'0.68%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001296"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'2.65%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.005912"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-3.49%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'3.403"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.88%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'4.445"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.23%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.3319"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.09%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'8.945"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-12.83%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1353"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-2.59%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.2552"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-1.11%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04159"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'1.16%"
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'-1.10%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004427"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'1.75%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001304"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'1.61%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0002987"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-20.32%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02731"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'3.87%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05621"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'0.39%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.007837"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.93%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1414"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.99%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.007373"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.27%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002115"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'10.71%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.007905"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-9.33%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3511"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006928"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-2.67%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000752"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.03%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003504"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-0.47%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003539"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'40.12%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002106"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.03%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002005"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.03%"
$ws.Range("E51").Style = "Normal"
